$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.106.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.790.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5184"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07993"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.266"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.788.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001091"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06556"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.942"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.135.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.271"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.993.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.330"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "122.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1082"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.670"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.525"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07192"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.65%  "
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.559"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("E41").Value = "  -2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.161"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.364"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.760"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5927"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.214"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.914"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06754"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.55%  "
